# Apply edits described by the diff to ValueSet-StrokeToiletingVS.xlsx

$wb = $excel.ActiveWorkbook

# --- Rename the "Include from " sheet ---
$wsInclude = $wb.Worksheets.Item("Include from ")
$wsInclude.Name = "Include from Prestroke and Po"

# --- Metadata sheet updates ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# Version: 1.0.1 -> 0.0.0
$wsMeta.Range("B3").Value = "0.0.0"

# Title: "ValueSet of Prestroke and Poststroke Functional Status: Toileting"
#        -> "Prestroke and Poststroke Functional Status: Toileting"
$wsMeta.Range("B5").Value = "Prestroke and Poststroke Functional Status: Toileting"

# Experimental: blank -> "false" (stored as literal text, not boolean).
# Assigning the literal word "false" straight to .Value auto-converts the
# cell to a Boolean, so the text is staged on a scratch cell (forced to
# stay text with a leading apostrophe) and copied in as values-only; that
# keeps B7's original number format/style untouched.
$scratch = $wsMeta.Range("Z1")
$scratch.Value = "'false"
$scratch.Copy()
$wsMeta.Range("B7").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
$excel.CutCopyMode = $false

# Date: 2023-11-21T19:08:35-03:00 -> 2024-01-11T13:00:00-03:00
$wsMeta.Range("B8").Value = "2024-01-11T13:00:00-03:00"

# Description
$wsMeta.Range("B12").Value = "ValueSet that defines the response values for the Prestroke and Poststroke Functional Status: Toileting."

# --- "Include from Prestroke and Po" sheet updates ---
# System URI value
$wsInclude.Range("B5").Value = "https://molic-avc.gabriellesantosleandro.com/CodeSystem/StrokeFuncStatusCS"
